$wb = $excel.ActiveWorkbook

# --- Sheet "Nädal 5": add new time-log entries (rows 14 and 15) ---
$ws5 = $wb.Worksheets.Item("Nädal 5")

# Row 14: fill in Stop time, Delta Time (minutes) and Aids ("x")
$ws5.Range("D14").Value = 0.48958333333333331
$ws5.Range("F14").Value = 195
$ws5.Range("I14").Value = "x"

# Update the wording of the existing comment in row 14 (shared string edit)
$ws5.Range("H14").Value = "kodutöö kaitsmine, täiendamine, v10"

# Row 15: fill in Start/Stop times, Delta Time, Activity and Comments
$ws5.Range("C15").Value = 0.93055555555555547
$ws5.Range("D15").Value = 0.99305555555555547
$ws5.Range("F15").Value = 90
$ws5.Range("G15").Value = "video"
$ws5.Range("H15").Value = "V11"

# --- Update sheetView selections on the various sheets (recorded cursor moves) ---
$ws3 = $wb.Worksheets.Item("Nädal 3")
$ws3.Activate() | Out-Null
$ws3.Range("G18").Select() | Out-Null

$ws4 = $wb.Worksheets.Item("Nädal 4")
$ws4.Activate() | Out-Null
$ws4.Range("F20").Select() | Out-Null

$ws5.Activate() | Out-Null
$ws5.Range("I15").Select() | Out-Null
